$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('salexandru', '3.10.2020, 16:57:24', '3.10.2020, 16:57:26'),
    @('jane01', '3.10.2020, 20:8:40', '3.10.2020, 20:9:43'),
    @('jane01', '3.10.2020, 20:23:30', '3.10.2020, 20:23:44'),
    @('jane01', '3.10.2020, 20:23:50', '3.10.2020, 20:23:53'),
    @('jane01', '3.10.2020, 20:52:52', '3.10.2020, 20:56:53'),
    @('jane01', '3.10.2020, 20:57:48', '3.10.2020, 20:57:51'),
    @('jane01', '3.10.2020, 20:57:53', '3.10.2020, 20:58:4'),
    @('jane01', '3.10.2020, 20:58:9', '3.10.2020, 20:58:29'),
    @('jane01', '3.10.2020, 21:9:32', '3.10.2020, 21:9:49'),
    @('jane01', '3.10.2020, 21:10:8', '3.10.2020, 21:10:10'),
    @('jane01', '3.10.2020, 21:11:4', '3.10.2020, 21:11:6'),
    @('jane01', '3.10.2020, 21:11:25', '3.10.2020, 21:11:26'),
    @('jane01', '3.10.2020, 21:11:36', '3.10.2020, 21:11:37'),
    @('jane01', '3.10.2020, 21:14:1', '3.10.2020, 21:15:4'),
    @('jane01', '3.10.2020, 21:16:8', '3.10.2020, 21:16:13'),
    @('jane01', '3.10.2020, 21:16:24', '3.10.2020, 21:16:34'),
    @('jane01', '3.10.2020, 21:16:37', '3.10.2020, 21:16:41'),
    @('jane01', '3.10.2020, 21:17:8', '3.10.2020, 21:17:12'),
    @('jane01', '3.10.2020, 21:18:16', '3.10.2020, 21:19:0'),
    @('jane01', '3.10.2020, 21:19:5', '3.10.2020, 21:20:10'),
    @('jane01', '3.10.2020, 21:20:20', '3.10.2020, 21:21:3'),
    @('jane01 (cu dwawdd)', '3.10.2020, 21:21:12', '3.10.2020, 21:23:27'),
    @('jane01 (cu ddd)', '3.10.2020, 21:23:30', '3.10.2020, 21:23:35'),
    @('jane01', '3.10.2020, 21:23:52', '3.10.2020, 21:23:57'),
    @('jane01', '3.10.2020, 21:27:1', '3.10.2020, 21:28:9'),
    @('salexandru (cu Alex)', '3.10.2020, 21:27:18', '3.10.2020, 21:27:28')
)

$startRow = 8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
